$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and B (rows 2-10) are formatted as Text (numFmtId 49) but the
# stored values are plain numbers. Temporarily switch to General so the
# COM layer doesn't coerce the new numbers into text, then restore the
# original Text format (matches the original file's cell styling).
$numRange = $ws.Range("A2:B10")
$numRange.NumberFormat = "General"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 6

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 8

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 9

$numRange.NumberFormat = "@"

# Update the Notes column (C) text values for rows 2-11
$ws.Range("C2").Value = "notes1"
$ws.Range("C3").Value = "notes2"
$ws.Range("C4").Value = "notes3"
$ws.Range("C5").Value = "notes4"
$ws.Range("C6").Value = "notes5"
$ws.Range("C7").Value = "notes6"
$ws.Range("C8").Value = "notes7"
$ws.Range("C9").Value = "notes8"
$ws.Range("C10").Value = "notes9"

# Row 11: A11/B11 become the text "10", C11 becomes "notes10"
$ws.Range("A11").Value = "10"
$ws.Range("B11").Value = "10"
$ws.Range("C11").Value = "notes10"

# Update the selection on the sheet to match the new edit
$ws.Range("A2:C11").Select()
